# Append a new data row (row 57) to each of the four sheets, carrying
# forward the same layout/format as the existing rows (A = datetime
# serial with the "YYYY-MM-DD HH:MM:SS" format, B-E = text byte-strings,
# F-I = plain numbers).

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = 45752.86331204861
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x6e"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = 45752.71608376157
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x6e"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = 45752.85688471065
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x6e"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 366
        I = 3
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = 45752.91387545139
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x6e"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 366
        I = 3
    }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)

    # Existing data occupies rows 1 (header) .. 56; the new record goes
    # immediately below the current last row.
    $lastRow = $ws.UsedRange.Rows.Count
    $newRow = $lastRow + 1

    $ws.Cells.Item($newRow, 1).Value2 = $r.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($newRow, 2).Value = $r.B
    $ws.Cells.Item($newRow, 3).Value = $r.C
    $ws.Cells.Item($newRow, 4).Value = $r.D
    $ws.Cells.Item($newRow, 5).Value = $r.E
    $ws.Cells.Item($newRow, 6).Value = $r.F
    $ws.Cells.Item($newRow, 7).Value = $r.G
    $ws.Cells.Item($newRow, 8).Value = $r.H
    $ws.Cells.Item($newRow, 9).Value = $r.I
}
